$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-73: I = value, J = value
$iVals = @(
    5,
    7,
    6,
    8,
    11,
    8,
    7,
    9,
    7,
    10,
    6,
    8,
    7,
    9,
    4,
    7,
    6,
    7,
    9,
    11,
    10,
    8,
    8,
    6,
    5,
    12,
    9,
    5,
    3,
    4,
    6,
    8,
    8,
    7,
    3,
    5,
    8,
    8,
    6,
    9,
    9,
    7,
    7,
    7,
    6,
    7,
    7,
    2,
    9,
    6,
    6,
    5,
    8,
    10,
    7,
    5,
    8,
    5,
    5,
    6,
    6,
    9,
    5,
    6,
    6,
    7,
    7,
    8,
    3,
    8,
    4,
    5
)
$jVals = @(
    6,
    7,
    6,
    8,
    11,
    8,
    8,
    9,
    7,
    10,
    7,
    8,
    7,
    9,
    5,
    7,
    6,
    8,
    9,
    11,
    11,
    8,
    9,
    7,
    5,
    12,
    9,
    5,
    3,
    4,
    6,
    9,
    8,
    7,
    3,
    6,
    8,
    8,
    6,
    9,
    9,
    7,
    7,
    7,
    6,
    8,
    7,
    3,
    9,
    7,
    6,
    5,
    8,
    10,
    7,
    6,
    9,
    6,
    5,
    6,
    7,
    9,
    5,
    6,
    6,
    7,
    7,
    8,
    4,
    8,
    4,
    5
)

for ($r = 2; $r -le 73; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

# Header row: add I0 and IF labels, matching the style of the existing header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
